$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.866.81"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.736.69"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.34"
$ws.Range("E5").Value = "  +5.19%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5206"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2740"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06156"
$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.739.65"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07178"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.00"
$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6430"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.611"
$ws.Range("E14").Value = "  +1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.20"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.903.23"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.76"
$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006772"
$ws.Range("E20").Value = "  +1.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.962.05"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.280"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.626"
$ws.Range("E23").Value = "  -1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.273"
$ws.Range("E24").Value = "  +2.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.11"
$ws.Range("E25").Value = "  -2.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.513"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.19"
$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.771"
$ws.Range("E28").Value = "  -1.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.01"
$ws.Range("E29").Value = "  +2.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.945"
$ws.Range("E30").Value = "  +5.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08248"
$ws.Range("E31").Value = "  -0.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.657"
$ws.Range("E32").Value = "  +3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04661"
$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.661"
$ws.Range("E34").Value = "  +1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9901"
$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6190"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01598"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.923"
$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9994"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.45"
$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3858"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7467"
$ws.Range("E43").Value = "  +1.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.008"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1124"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.261"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47: Cronos -> Aave
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.07"
$ws.Range("E47").Value = "  +2.88%  "

# Row 48: Aave -> Cronos
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05217"
$ws.Range("E48").Value = "  -2.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.60"
$ws.Range("E49").Value = "  +1.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.577"
$ws.Range("E50").Value = "  -1.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3414"
$ws.Range("E51").Value = "  -0.20%  "
